$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at position 4 (new FA "losfahren" requirement), pushing
#    the existing FA2..NFA2 block (and the trailing blank row) down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# Copy the formatting of row 3 (style pair used by plain FA rows) onto the
# freshly inserted row 4.
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4 itself stays label-less; only column B gets the new requirement text.
$ws.Range("B4").Value = "Die Fahrzeuge sollen dem Server Informationen senden, wenn sie losfahren."

# ---------------------------------------------------------------------------
# 2) Drop the stray note that used to live in column G (old row 7, now row 8
#    after the insert above).
# ---------------------------------------------------------------------------
$ws.Range("G8").ClearContents()

# ---------------------------------------------------------------------------
# 3) Insert three new rows above the old "NFA2 / Python" row (now row 14)
#    to make room for a reworded NFA2, a brand-new NFA3, and a spacer row.
# ---------------------------------------------------------------------------
$ws.Range("A14:A16").EntireRow.Insert()

# Rows 14 and 15 get the "sub-item" style pair (bold-ish A column + plain B).
$ws.Range("A6:B6").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)
$ws.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 16 is a plain blank spacer row, matching row 2's style.
$ws.Range("A2:B2").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A14").Value = "NFA2"
$ws.Range("B14").Value = "Das Versenden von Informationen sollte nicht länger wie 5s dauern."

$ws.Range("A15").Value = "NFA3"
$ws.Range("B15").Value = "Das erhalten von Informationen sollte nicht länger wie 5s dauern."

# Row 17 (old "NFA2 / Python" row, unchanged formatting) becomes NFA4.
$ws.Range("A17").Value = "NFA4"

# Leave the selection where the author left it when they saved the file.
$ws.Range("B29").Select() | Out-Null
